$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
# B2: text changed in place (keep existing formatting)
$ws.Range("B2").Value = "report_master"

# C2 / D2: two new values added to the right of B2
$ws.Range("C2").Value = "data"
$ws.Range("D2").Value = "report_id-25"

# E2: old value removed
$ws.Range("E2").ClearContents()

# --- Row 3 ---
$ws.Range("A3").Value = "kmb_reports"
$ws.Range("B3").Value = "report_template"
$ws.Range("C3").Value = "data"
$ws.Range("D3").Value = "report_matser-25"

# --- Selection moves to C4 ---
$ws.Range("C4").Select() | Out-Null
